$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2024-01-13 Saturday", $true, $false, $false, $false,
                         $false, $true, 1, $false, "2024-01-14 Sunday", 2)

# Update the division problems in the table, cell by cell (row, column),
# since several expressions repeat and a global replace would be ambiguous.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text  = "96÷2="
$t.Cell(1,2).Range.Text  = "52÷2="
$t.Cell(1,3).Range.Text  = "44÷7="
$t.Cell(1,4).Range.Text  = "99÷8="
$t.Cell(1,5).Range.Text  = "31÷6="

$t.Cell(5,1).Range.Text  = "30÷5="
$t.Cell(5,2).Range.Text  = "81÷6="
$t.Cell(5,3).Range.Text  = "84÷4="
$t.Cell(5,4).Range.Text  = "33÷6="
$t.Cell(5,5).Range.Text  = "70÷5="

$t.Cell(9,1).Range.Text  = "58÷3="
$t.Cell(9,2).Range.Text  = "53÷4="
$t.Cell(9,3).Range.Text  = "37÷7="
$t.Cell(9,4).Range.Text  = "46÷8="
$t.Cell(9,5).Range.Text  = "30÷2="

$t.Cell(13,1).Range.Text = "46÷6="
$t.Cell(13,2).Range.Text = "78÷4="
$t.Cell(13,3).Range.Text = "14÷6="
$t.Cell(13,4).Range.Text = "64÷3="
$t.Cell(13,5).Range.Text = "52÷3="

$t.Cell(17,1).Range.Text = "44÷5="
$t.Cell(17,2).Range.Text = "64÷7="
$t.Cell(17,3).Range.Text = "63÷3="
$t.Cell(17,4).Range.Text = "18÷6="
$t.Cell(17,5).Range.Text = "99÷5="
